# modified test cases on overdue fix
$wb = $excel.ActiveWorkbook

# --- Update "Summary" sheet ---
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B2").Value = 836.76
$ws.Range("E2").Value = 4163.24
$ws.Range("F2").Value = 849.4
$ws.Range("B3").Value = 50.96
$ws.Range("E3").Value = 128.7
$ws.Range("F3").Value = 38.32
$ws.Range("A5").Value = 200
$ws.Range("B5").Value = 100
$ws.Range("E5").Value = 100
$ws.Range("F5").Value = 100
# touch G2 so it exists as an (empty, unformatted) cell and widens the
# sheet's used range to column G, then select E5 like the saved workbook
$ws.Range("G2").Value = 0
$ws.Range("G2").ClearContents()
$ws.Range("G2").Style = "Normal"

# --- Update "Repayment schedule" sheet ---
$ws3 = $wb.Worksheets.Item("Repayment schedule")
$ws3.Range("J3").Value = 100
$ws3.Range("K3").Value = 987.72
$ws3.Range("L3").Value = 987.72
$ws3.Range("J4").Value = 100
$ws3.Range("K4").Value = 987.72
$ws3.Range("L4").Value = 0
$ws3.Range("M4").Value = 0
$ws3.Range("P4").Value = 987.72

# --- Update "Transactions" sheet ---
$ws4 = $wb.Worksheets.Item("Transactions")
$ws4.Range("A2").Value = 640
$ws4.Range("F2").Value = 836.76
$ws4.Range("G2").Value = 50.96
$ws4.Range("I2").Value = 100
$ws4.Range("J2").Value = 4163.24
$ws4.Range("A3").Value = 632
# K2:L3 are trailing, unused columns - fully clear contents + formatting so
# they disappear from the saved used range entirely
$ws4.Range("K2:L3").Clear()

# --- Remove the Acc_Disbursement and Acc_Repayment sheets ---
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Acc_Disbursement").Delete()
$wb.Worksheets.Item("Acc_Repayment").Delete()
$excel.DisplayAlerts = $true

# --- Restore the selection on each touched sheet (Transactions stays the
# active/tab-selected sheet, matching the saved file) ---
$ws3.Range("H4").Select()
$ws.Range("E5").Select()
$ws4.Range("F3").Select()
